# Add a new "date block" (45901) to the bottom of the data table, mirroring
# the structure of the previous block (rows 170:190, date 45870), and add
# the blank J/K/L cells (columns 10-12) to that previous block.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------
# 1) Duplicate rows 170:190 (the 45870 block) down to 191:211 so that the
#    B/C/D (asset/benchmark/description) text and A/E:I formatting carries
#    over exactly as in the source block.
# ---------------------------------------------------------------------
$src = $ws.Range("A170:I190")
$dst = $ws.Range("A191:I211")
$src.Copy($dst)

# ---------------------------------------------------------------------
# 2) Overwrite column A (date) and E:I (weight/perf figures) for the new
#    block with the 2025-09 values.
# ---------------------------------------------------------------------
$newDate = 45901

$values = @{
  191 = @("0.11019231603227203","3.319084385822061E-4","3.3291030278731482E-4","3.2933728454145302E-3","-4.4979670100069834E-3")
  192 = @("0.47525833850012411","4.9902125542513037E-3","6.5442621134771707E-3","2.9846223657807792E-2","4.601307848159715E-2")
  193 = @("0.19354034966519451","2.187005951216698E-3","2.3774772357617081E-3","1.1980147644275539E-2","1.7755794959571709E-2")
  194 = @("4.393808970867856E-2","3.9104899840723917E-4","3.874624558066036E-4","3.1986929307917995E-3","3.1053064816854974E-3")
  195 = @("0","0","0","0","-1.1914618342282394E-3")
  196 = @("4.4943551969111727E-2","7.0561376591505404E-4","1.4575562874360689E-3","4.8314318366795101E-3","7.2827138846488851E-3")
  197 = @("0","0","0","0","3.3527321496243795E-3")
  198 = @("9.8989519553941183E-2","9.997941474948059E-4","1.3899521979822476E-3","6.6669441419579378E-3","1.0498166205078069E-2")
  199 = @("0","0","0","0","0")
  200 = @("9.3846827603198124E-2","9.5723764155262092E-4","8.9999999999999998E-4","5.0583440078123794E-3","5.1999999999999998E-3")
  201 = @("0.226728484630008","8.0715340528282834E-3","1.0436093170480894E-2","3.3238395846759174E-2","3.0963135377441484E-2")
  202 = @("4.2949175287504485E-2","1.5289906402351596E-3","1.0154406222492782E-3","6.2963490971481581E-3","1.4335303688439889E-3")
  203 = @("0.12961474215734423","4.6142848208014543E-3","5.0350026728010304E-3","1.9001521200266665E-2","1.9031764051594409E-2")
  204 = @("0","0","0","0","0")
  205 = @("2.7124124082333268E-2","9.6561881733106438E-4","1.3680447121012297E-3","3.9763965904700575E-3","5.1656012810432517E-3")
  206 = @("2.7040443102826012E-2","9.6263977446060605E-4","3.0176051633293544E-3","3.9641289588742941E-3","5.3322396759598305E-3")
  207 = @("0.10704515670200816","2.0231534616679545E-3","-1.972687351531666E-4","7.8785235332677998E-3","7.6E-3")
  208 = @("5.8366193223647833E-2","1.0447548587032961E-3","5.0000000000000001E-3","3.4669518774846815E-3","1.6199999999999999E-2")
  209 = @("2.2192661670125881E-2","9.9866977515566453E-5","1E-4","1.3670679588797543E-3","1.39E-3")
  210 = @("0","0","0","0","0")
  211 = @("2.0041976270155602E-2","3.5474297998175418E-4","1E-4","1.9581010815942022E-3","9.0000000000000006E-5")
}

for ($r = 191; $r -le 211; $r++) {
  $ws.Cells.Item($r, 1).Value2 = $newDate
  $row = $values[$r]
  $ws.Cells.Item($r, 5).Value2 = [double]$row[0]
  $ws.Cells.Item($r, 6).Value2 = [double]$row[1]
  $ws.Cells.Item($r, 7).Value2 = [double]$row[2]
  $ws.Cells.Item($r, 8).Value2 = [double]$row[3]
  $ws.Cells.Item($r, 9).Value2 = [double]$row[4]
}

# ---------------------------------------------------------------------
# 3) Add the (empty) J/K/L columns to the previous block (rows 170:190),
#    matching the styles used elsewhere in the sheet (percentage format,
#    with J/L highlighted blue like column G/I, K plain like F/H).
# ---------------------------------------------------------------------
$ws.Range("E2").Copy()
$ws.Range("J170:J190").PasteSpecial(-4122)
$ws.Range("L170:L190").PasteSpecial(-4122)

$ws.Range("F2").Copy()
$ws.Range("K170:K190").PasteSpecial(-4122)

# ---------------------------------------------------------------------
# 4) Cosmetic touch-ups matching the author's final view state: widen
#    column A (now that it holds wider content) and move the active
#    selection / scroll position.
# ---------------------------------------------------------------------
$ws.Columns.Item(1).ColumnWidth = 9.5

$ws.Activate()
$ws.Range("D204").Select()
